# Resolve edit employee issue
# - Update existing employee E0111 (row 5) Current Role to "Solution Developer"
# - Add 6 new employee rows (7-12) to the Employees sheet, copying the formatting
#   from the last existing row (row 6)
# - Move the active selection to F12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# Fix Current Role for existing employee in row 5 (E0111 / p11)
$ws.Range("F5").Value = "Solution Developer"

# Seed rows 7-12 by copying row 6's formatting, then overwrite values
$ws.Range("A6:K6").Copy($ws.Range("A7:K7"))
$ws.Range("A6:K6").Copy($ws.Range("A8:K8"))
$ws.Range("A6:K6").Copy($ws.Range("A9:K9"))
$ws.Range("A6:K6").Copy($ws.Range("A10:K10"))
$ws.Range("A6:K6").Copy($ws.Range("A11:K11"))
$ws.Range("A6:K6").Copy($ws.Range("A12:K12"))

# Row 7: E0115 / p15
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "E0115"
$ws.Range("C7").Value = "p15"
$ws.Range("D7").Value = 44922
$ws.Range("E7").Value = 35954
$ws.Range("F7").Value = "Java Developer"
$ws.Range("G7").Value = "p15@gmail.com"
$ws.Range("H7").Value = "Female"
$ws.Range("I7").Value = 1234543268
$ws.Range("J7").Value = "Pune"
$ws.Range("K7").Value = 2

# Row 8: E0119 / p19
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "E0119"
$ws.Range("C8").Value = "p19"
$ws.Range("D8").Value = 44922
$ws.Range("E8").Value = 35967
$ws.Range("F8").Value = "Java Developer"
$ws.Range("G8").Value = "p19@gmail.com"
$ws.Range("H8").Value = "Female"
$ws.Range("I8").Value = 1234543268
$ws.Range("J8").Value = "Pune"
$ws.Range("K8").Value = 1

# Row 9: E0120 / p20
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "E0120"
$ws.Range("C9").Value = "p20"
$ws.Range("D9").Value = 44922
$ws.Range("E9").Value = 35959
$ws.Range("F9").Value = "Solution Developer"
$ws.Range("G9").Value = "p20@gmail.com"
$ws.Range("H9").Value = "Female"
$ws.Range("I9").Value = 1234543268
$ws.Range("J9").Value = "Pune"
$ws.Range("K9").Value = 1

# Row 10: E0121 / p21
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "E0121"
$ws.Range("C10").Value = "p21"
$ws.Range("D10").Value = 44922
$ws.Range("E10").Value = 35964
$ws.Range("F10").Value = "Java Developer"
$ws.Range("G10").Value = "p21@gmail.com"
$ws.Range("H10").Value = "Female"
$ws.Range("I10").Value = 1234543268
$ws.Range("J10").Value = "Pune"
$ws.Range("K10").Value = 1

# Row 11: E0122 / p22 (Mobile left blank)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "E0122"
$ws.Range("C11").Value = "p22"
$ws.Range("D11").Value = 44922
$ws.Range("E11").Value = 35964
$ws.Range("F11").Value = "Solution Developer"
$ws.Range("G11").Value = "p22@gmail.com"
$ws.Range("H11").Value = "Female"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "Pune"
$ws.Range("K11").Value = 1

# Row 12: E0123 / p23 (Current Role, Mobile, Location left blank)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "E0123"
$ws.Range("C12").Value = "p23"
$ws.Range("D12").Value = 44922
$ws.Range("E12").Value = 35964
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = "p23@gmail.com"
$ws.Range("H12").Value = "Female"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = 1

# Move selection to match author's final cursor position
$ws.Range("F12").Select() | Out-Null
